$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.258.59"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.152.81"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.83%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.146.80"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").Value = "3.666.87"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "3.143.83"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "63.186.29"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.28%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.92%  "
$ws.Range("E34").Value = "  -5.69%  "
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").Value = "0.0₃0700"
$ws.Range("E38").Value = "  -8.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0391"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "2.928.96"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("E44").Value = "  -5.62%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
